$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 43.34730933333333
$ws.Range("H2").Value = 130.041928
$ws.Range("I2").Value = 0.04273139820300816
$ws.Range("J2").Value = 0.04273139820300816
$ws.Range("M2").Value = 0.6574793333333333
$ws.Range("N2").Value = 1.972438
$ws.Range("O2").Value = 0.04234443143670402
$ws.Range("P2").Value = 0.04234443143670403
$ws.Range("Q2").Value = 28.49996004227377
$ws.Range("R2").Value = 256.499640380464
$ws.Range("S2").Value = 0.001809436761401777
$ws.Range("T2").Value = 0.001809436761401777

$ws.Range("G3").Value = 43.34730933333333
$ws.Range("H3").Value = 130.041928
$ws.Range("I3").Value = 0.04273139820300816
$ws.Range("J3").Value = 0.04273139820300816
$ws.Range("M3").Value = 1.851391
$ws.Range("N3").Value = 5.554173
$ws.Range("O3").Value = 0.1192373589365509
$ws.Range("P3").Value = 0.119237358936551
$ws.Range("Q3").Value = 80.25281837394931
$ws.Range("R3").Value = 722.2753653655438
$ws.Range("S3").Value = 0.005095179065392773
$ws.Range("T3").Value = 0.005095179065392772

$ws.Range("G4").Value = 43.34730933333333
$ws.Range("H4").Value = 130.041928
$ws.Range("I4").Value = 0.04273139820300816
$ws.Range("J4").Value = 0.04273139820300816
$ws.Range("M4").Value = 5.370269333333333
$ws.Range("N4").Value = 16.110808
$ws.Range("O4").Value = 0.3458679080132824
$ws.Range("P4").Value = 0.3458679080132824
$ws.Range("Q4").Value = 232.7867259953137
$ws.Range("R4").Value = 2095.080533957824
$ws.Range("S4").Value = 0.01477941930295697
$ws.Range("T4").Value = 0.01477941930295697

$ws.Range("G5").Value = 43.34730933333333
$ws.Range("H5").Value = 130.041928
$ws.Range("I5").Value = 0.04273139820300816
$ws.Range("J5").Value = 0.04273139820300816
$ws.Range("M5").Value = 1.801189666666667
$ws.Range("N5").Value = 5.403569
$ws.Range("O5").Value = 0.1160041821512257
$ws.Range("P5").Value = 0.1160041821512257
$ws.Range("Q5").Value = 78.07672564900355
$ws.Range("R5").Value = 702.6905308410319
$ws.Range("S5").Value = 0.004957020900718318
$ws.Range("T5").Value = 0.004957020900718318

$ws.Range("G6").Value = 43.34730933333333
$ws.Range("H6").Value = 130.041928
$ws.Range("I6").Value = 0.04273139820300816
$ws.Range("J6").Value = 0.04273139820300816
$ws.Range("M6").Value = 5.846608
$ws.Range("N6").Value = 17.539824
$ws.Range("O6").Value = 0.3765461194622369
$ws.Range("P6").Value = 0.376546119462237
$ws.Range("Q6").Value = 253.4347255267413
$ws.Range("R6").Value = 2280.912529740672
$ws.Range("S6").Value = 0.01609034217253833
$ws.Range("T6").Value = 0.01609034217253833

$ws.Range("G7").Value = 90.21844233333333
$ws.Range("H7").Value = 270.655327
$ws.Range("I7").Value = 0.0889365509391893
$ws.Range("J7").Value = 0.08893655093918929
$ws.Range("M7").Value = 0.6574793333333333
$ws.Range("N7").Value = 1.972438
$ws.Range("O7").Value = 0.04234443143670402
$ws.Range("P7").Value = 0.04234443143670403
$ws.Range("Q7").Value = 59.31676131969177
$ws.Range("R7").Value = 533.850851877226
$ws.Range("S7").Value = 0.003765967683461436
$ws.Range("T7").Value = 0.003765967683461436

$ws.Range("G8").Value = 90.21844233333333
$ws.Range("H8").Value = 270.655327
$ws.Range("I8").Value = 0.0889365509391893
$ws.Range("J8").Value = 0.08893655093918929
$ws.Range("M8").Value = 1.851391
$ws.Range("N8").Value = 5.554173
$ws.Range("O8").Value = 0.1192373589365509
$ws.Range("P8").Value = 0.119237358936551
$ws.Range("Q8").Value = 167.0296121699523
$ws.Range("R8").Value = 1503.266509529571
$ws.Range("S8").Value = 0.01060455944691496
$ws.Range("T8").Value = 0.01060455944691496

$ws.Range("G9").Value = 90.21844233333333
$ws.Range("H9").Value = 270.655327
$ws.Range("I9").Value = 0.0889365509391893
$ws.Range("J9").Value = 0.08893655093918929
$ws.Range("M9").Value = 5.370269333333333
$ws.Range("N9").Value = 16.110808
$ws.Range("O9").Value = 0.3458679080132824
$ws.Range("P9").Value = 0.3458679080132824
$ws.Range("Q9").Value = 484.4973341638017
$ws.Range("R9").Value = 4360.476007474215
$ws.Range("S9").Value = 0.03076029881925413
$ws.Range("T9").Value = 0.03076029881925412

$ws.Range("G10").Value = 90.21844233333333
$ws.Range("H10").Value = 270.655327
$ws.Range("I10").Value = 0.0889365509391893
$ws.Range("J10").Value = 0.08893655093918929
$ws.Range("M10").Value = 1.801189666666667
$ws.Range("N10").Value = 5.403569
$ws.Range("O10").Value = 0.1160041821512257
$ws.Range("P10").Value = 0.1160041821512257
$ws.Range("Q10").Value = 162.5005260735626
$ws.Range("R10").Value = 1462.504734662063
$ws.Range("S10").Value = 0.01031701185505148
$ws.Range("T10").Value = 0.01031701185505148

$ws.Range("G11").Value = 90.21844233333333
$ws.Range("H11").Value = 270.655327
$ws.Range("I11").Value = 0.0889365509391893
$ws.Range("J11").Value = 0.08893655093918929
$ws.Range("M11").Value = 5.846608
$ws.Range("N11").Value = 17.539824
$ws.Range("O11").Value = 0.3765461194622369
$ws.Range("P11").Value = 0.376546119462237
$ws.Range("Q11").Value = 527.4718666936053
$ws.Range("R11").Value = 4747.246800242448
$ws.Range("S11").Value = 0.03348871313450729
$ws.Range("T11").Value = 0.03348871313450729

$ws.Range("G12").Value = 394.701121
$ws.Range("H12").Value = 1184.103363
$ws.Range("I12").Value = 0.3890929110023202
$ws.Range("J12").Value = 0.3890929110023201
$ws.Range("M12").Value = 0.6574793333333333
$ws.Range("N12").Value = 1.972438
$ws.Range("O12").Value = 0.04234443143670402
$ws.Range("P12").Value = 0.04234443143670403
$ws.Range("Q12").Value = 259.5078299009993
$ws.Range("R12").Value = 2335.570469108994
$ws.Range("S12").Value = 0.01647591809244533
$ws.Range("T12").Value = 0.01647591809244532

$ws.Range("G13").Value = 394.701121
$ws.Range("H13").Value = 1184.103363
$ws.Range("I13").Value = 0.3890929110023202
$ws.Range("J13").Value = 0.3890929110023201
$ws.Range("M13").Value = 1.851391
$ws.Range("N13").Value = 5.554173
$ws.Range("O13").Value = 0.1192373589365509
$ws.Range("P13").Value = 0.119237358936551
$ws.Range("Q13").Value = 730.7461031093109
$ws.Range("R13").Value = 6576.714927983799
$ws.Range("S13").Value = 0.04639441108885112
$ws.Range("T13").Value = 0.04639441108885112

$ws.Range("G14").Value = 394.701121
$ws.Range("H14").Value = 1184.103363
$ws.Range("I14").Value = 0.3890929110023202
$ws.Range("J14").Value = 0.3890929110023201
$ws.Range("M14").Value = 5.370269333333333
$ws.Range("N14").Value = 16.110808
$ws.Range("O14").Value = 0.3458679080132824
$ws.Range("P14").Value = 0.3458679080132824
$ws.Range("Q14").Value = 2119.651325938589
$ws.Range("R14").Value = 19076.8619334473
$ws.Range("S14").Value = 0.1345747511511707
$ws.Range("T14").Value = 0.1345747511511707

$ws.Range("G15").Value = 394.701121
$ws.Range("H15").Value = 1184.103363
$ws.Range("I15").Value = 0.3890929110023202
$ws.Range("J15").Value = 0.3890929110023201
$ws.Range("M15").Value = 1.801189666666667
$ws.Range("N15").Value = 5.403569
$ws.Range("O15").Value = 0.1160041821512257
$ws.Range("P15").Value = 0.1160041821512257
$ws.Range("Q15").Value = 710.9315805669497
$ws.Range("R15").Value = 6398.384225102547
$ws.Range("S15").Value = 0.04513640492166381
$ws.Range("T15").Value = 0.0451364049216638

$ws.Range("G16").Value = 394.701121
$ws.Range("H16").Value = 1184.103363
$ws.Range("I16").Value = 0.3890929110023202
$ws.Range("J16").Value = 0.3890929110023201
$ws.Range("M16").Value = 5.846608
$ws.Range("N16").Value = 17.539824
$ws.Range("O16").Value = 0.3765461194622369
$ws.Range("P16").Value = 0.376546119462237
$ws.Range("Q16").Value = 2307.662731647568
$ws.Range("R16").Value = 20768.96458482811
$ws.Range("S16").Value = 0.1465114257481892
$ws.Range("T16").Value = 0.1465114257481892

$ws.Range("G17").Value = 7.804371333333333
$ws.Range("H17").Value = 23.413114
$ws.Range("I17").Value = 0.007693480963358413
$ws.Range("J17").Value = 0.007693480963358412
$ws.Range("M17").Value = 0.6574793333333333
$ws.Range("N17").Value = 1.972438
$ws.Range("O17").Value = 0.04234443143670402
$ws.Range("P17").Value = 0.04234443143670403
$ws.Range("Q17").Value = 5.131212861325777
$ws.Range("R17").Value = 46.180915751932
$ws.Range("S17").Value = 0.0003257760771625179
$ws.Range("T17").Value = 0.0003257760771625179

$ws.Range("G18").Value = 7.804371333333333
$ws.Range("H18").Value = 23.413114
$ws.Range("I18").Value = 0.007693480963358413
$ws.Range("J18").Value = 0.007693480963358412
$ws.Range("M18").Value = 1.851391
$ws.Range("N18").Value = 5.554173
$ws.Range("O18").Value = 0.1192373589365509
$ws.Range("P18").Value = 0.119237358936551
$ws.Range("Q18").Value = 14.44894284719133
$ws.Range("R18").Value = 130.040485624722
$ws.Range("S18").Value = 0.0009173503510994888
$ws.Range("T18").Value = 0.0009173503510994888

$ws.Range("G19").Value = 7.804371333333333
$ws.Range("H19").Value = 23.413114
$ws.Range("I19").Value = 0.007693480963358413
$ws.Range("J19").Value = 0.007693480963358412
$ws.Range("M19").Value = 5.370269333333333
$ws.Range("N19").Value = 16.110808
$ws.Range("O19").Value = 0.3458679080132824
$ws.Range("P19").Value = 0.3458679080132824
$ws.Range("Q19").Value = 41.91157603734577
$ws.Range("R19").Value = 377.204184336112
$ws.Range("S19").Value = 0.002660928166136787
$ws.Range("T19").Value = 0.002660928166136787

$ws.Range("G20").Value = 7.804371333333333
$ws.Range("H20").Value = 23.413114
$ws.Range("I20").Value = 0.007693480963358413
$ws.Range("J20").Value = 0.007693480963358412
$ws.Range("M20").Value = 1.801189666666667
$ws.Range("N20").Value = 5.403569
$ws.Range("O20").Value = 0.1160041821512257
$ws.Range("P20").Value = 0.1160041821512257
$ws.Range("Q20").Value = 14.05715300042955
$ws.Range("R20").Value = 126.514377003866
$ws.Range("S20").Value = 0.0008924759670504167
$ws.Range("T20").Value = 0.0008924759670504168

$ws.Range("G21").Value = 7.804371333333333
$ws.Range("H21").Value = 23.413114
$ws.Range("I21").Value = 0.007693480963358413
$ws.Range("J21").Value = 0.007693480963358412
$ws.Range("M21").Value = 5.846608
$ws.Range("N21").Value = 17.539824
$ws.Range("O21").Value = 0.3765461194622369
$ws.Range("P21").Value = 0.376546119462237
$ws.Range("Q21").Value = 45.62909987243733
$ws.Range("R21").Value = 410.661898851936
$ws.Range("S21").Value = 0.002896950401909203
$ws.Range("T21").Value = 0.002896950401909203

$ws.Range("G22").Value = 478.3423056666667
$ws.Range("H22").Value = 1435.026917
$ws.Range("I22").Value = 0.4715456588921241
$ws.Range("J22").Value = 0.471545658892124
$ws.Range("M22").Value = 0.6574793333333333
$ws.Range("N22").Value = 1.972438
$ws.Range("O22").Value = 0.04234443143670402
$ws.Range("P22").Value = 0.04234443143670403
$ws.Range("Q22").Value = 314.5001802348496
$ws.Range("R22").Value = 2830.501622113646
$ws.Range("S22").Value = 0.01996733282223297
$ws.Range("T22").Value = 0.01996733282223297

$ws.Range("G23").Value = 478.3423056666667
$ws.Range("H23").Value = 1435.026917
$ws.Range("I23").Value = 0.4715456588921241
$ws.Range("J23").Value = 0.471545658892124
$ws.Range("M23").Value = 1.851391
$ws.Range("N23").Value = 5.554173
$ws.Range("O23").Value = 0.1192373589365509
$ws.Range("P23").Value = 0.119237358936551
$ws.Range("Q23").Value = 885.5986396305157
$ws.Range("R23").Value = 7970.387756674641
$ws.Range("S23").Value = 0.05622585898429262
$ws.Range("T23").Value = 0.05622585898429261

$ws.Range("G24").Value = 478.3423056666667
$ws.Range("H24").Value = 1435.026917
$ws.Range("I24").Value = 0.4715456588921241
$ws.Range("J24").Value = 0.471545658892124
$ws.Range("M24").Value = 5.370269333333333
$ws.Range("N24").Value = 16.110808
$ws.Range("O24").Value = 0.3458679080132824
$ws.Range("P24").Value = 0.3458679080132824
$ws.Range("Q24").Value = 2568.82701495766
$ws.Range("R24").Value = 23119.44313461894
$ws.Range("S24").Value = 0.1630925105737638
$ws.Range("T24").Value = 0.1630925105737638

$ws.Range("G25").Value = 478.3423056666667
$ws.Range("H25").Value = 1435.026917
$ws.Range("I25").Value = 0.4715456588921241
$ws.Range("J25").Value = 0.471545658892124
$ws.Range("M25").Value = 1.801189666666667
$ws.Range("N25").Value = 5.403569
$ws.Range("O25").Value = 0.1160041821512257
$ws.Range("P25").Value = 0.1160041821512257
$ws.Range("Q25").Value = 861.5852180963083
$ws.Range("R25").Value = 7754.266962866774
$ws.Range("S25").Value = 0.05470126850674171
$ws.Range("T25").Value = 0.0547012685067417

$ws.Range("G26").Value = 478.3423056666667
$ws.Range("H26").Value = 1435.026917
$ws.Range("I26").Value = 0.4715456588921241
$ws.Range("J26").Value = 0.471545658892124
$ws.Range("M26").Value = 5.846608
$ws.Range("N26").Value = 17.539824
$ws.Range("O26").Value = 0.3765461194622369
$ws.Range("P26").Value = 0.376546119462237
$ws.Range("Q26").Value = 2796.679951049179
$ws.Range("R26").Value = 25170.11955944261
$ws.Range("S26").Value = 0.177558688005093
$ws.Range("T26").Value = 0.177558688005093
